# team_member.xlsx: rename the "电话" (Phone) column header to "分机号"
# (Extension number), and move the active selection to I10, matching the
# commit's OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("member")

# G1 currently holds the shared string "电话" -> change it to "分机号"
$ws.Range("G1").Value = "分机号"

# Reflect the sheetView's new active cell / selection (D6 -> I10)
$ws.Range("I10").Select()
